# qa.xlsx update — "support for 'How many' questions + fixes"
#
#  * B6 / B8 ("What is the population of Isle of Man / Djibouti?") switch
#    from plain numbers to comma-formatted text answers, matching the style
#    of the other answer cells ("1,499" etc.).
#  * The temporary yellow-highlight / scratch formatting that had been
#    applied to several rows (A3/B3, A19/B19, A21:B24, A27/B27) is cleared
#    back to the sheet's normal look.
#  * Column A is narrowed and no longer auto "best fit".
#  * Selection / scroll position left where the author ended up (C28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Turn the two population figures into comma-formatted text -------
# Setting NumberFormat to text ("@") first keeps Excel from re-parsing the
# comma-separated digits back into a number.
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value2 = "84,069"

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value2 = "921,804"

# --- 2. Clear the leftover highlight formatting on the affected rows ----
# Re-use the plain formatting already present elsewhere on the sheet
# (column A's normal cells use the same format as A4, column B's normal
# cells use the same format as B4, and B3 keeps the header-ish look of B1).
$ws.Range("A4").Copy() | Out-Null
foreach ($addr in @("A3", "A19", "A21", "A22", "A23", "A24", "A27")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

$ws.Range("B4").Copy() | Out-Null
foreach ($addr in @("B6", "B8", "B19", "B21", "B22", "B23", "B24", "B27")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

$ws.Range("B1").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null        # xlPasteFormats

$excel.CutCopyMode = 0

# --- 3. Resize column A, drop the "best fit" auto width ------------------
$ws.Columns.Item(1).ColumnWidth = 49.43

# --- 4. Leave the selection / scroll position where the author left it ---
$ws.Range("C28").Select() | Out-Null
$av = $excel.ActiveWindow
$av.ScrollRow = 21
$av.ScrollColumn = 1
